# "remove extra i from name"
#
# The document's default header (Sections(1).Headers(wdHeaderFooterPrimary))
# contains the line:
#       <tab>   Pavithra Viinay <tab>      Rittick Datta
# i.e. "Pavithra" followed by "Viinay" - the editor's name has an extra,
# duplicated "i" ("Vi" + "i" + "nay"). The fix removes the stray "i" so the
# name reads "Pavithra Vinay", and Word's "_GoBack" last-edit bookmark ends
# up sitting right after the fix (immediately after "...Pavithra V").

$d = $word.ActiveDocument
$hdr = $d.Sections(1).Headers(1)   # wdHeaderFooterPrimary

# Only touch the header if the "Viinay" typo is still there (keeps this
# script idempotent / safe to run more than once).
$checkRng = $hdr.Range
$hasTypo = $checkRng.Find.Execute("Viinay", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)

if ($hasTypo) {
    # --- Step 1: mark the edit point right after "Pavithra V" ---------
    # Collapsing Find's match range to its end gives us the exact
    # boundary between the "V" that should stay and the doubled "i"
    # that should go, and dropping a fresh "_GoBack" bookmark there
    # mirrors where Word leaves the last-edit marker after the fix.
    $markRng = $hdr.Range
    $foundMark = $markRng.Find.Execute("Pavithra V", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)

    if ($foundMark) {
        $editPoint = $markRng.End
        $bmRange = $hdr.Range.Duplicate
        $bmRange.Start = $editPoint
        $bmRange.End = $editPoint
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }

    # --- Step 2: collapse "Viinay" down to "Vinay" ---------------------
    # Replace the duplicated-i stretch (plus a touch of the following
    # text so the edit fully swallows the old "_GoBack" bookmark that
    # used to sit between the two i's) with the corrected spelling.
    $fixRng = $hdr.Range
    $fixRng.Find.Execute("iin", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "in", 2) | Out-Null
}
